# Updated symbol list on Fri Feb  3 06:44:27 UTC 2023 with GitHub Actions
# Refresh the Price (column D) and Volume(1h) (column E) values for the
# crypto rows on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as literal text (not auto-converted to a number
# or percentage) while keeping the cell's original (unformatted) style.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

# Rows whose Price (D) and Volume(1h) (E) both change
$priceAndVolume = @{
    2  = @{ D = "323.56";     E = "-2.02%" }
    3  = @{ D = "39.65";      E = "-1.26%" }
    4  = @{ D = "5.882";      E = "11.52%" }
    5  = @{ D = "0.08028";    E = "-0.85%" }
    6  = @{ D = "4.575";      E = "1.07%"  }
    7  = @{ D = "8.678";      E = "0.33%"  }
    8  = @{ D = "1.944";      E = "0.95%"  }
    10 = @{ D = "0.9296";     E = "-0.73%" }
    11 = @{ D = "0.1279";     E = "-4.43%" }
    12 = @{ D = "0.1970";     E = "0.11%"  }
    13 = @{ D = "8.714";      E = "34.06%" }
    14 = @{ D = "0.09124";    E = "0.34%"  }
    15 = @{ D = "0.03560";    E = "1.67%"  }
    16 = @{ D = "0.1047";     E = "9.18%"  }
    17 = @{ D = "0.001301";   E = "-7.30%" }
    18 = @{ D = "0.006141";   E = "-2.64%" }
    19 = @{ D = "3.349";      E = "-0.32%" }
    22 = @{ D = "0.2449";     E = "-4.79%" }
    23 = @{ D = "0.04411";    E = "-0.65%" }
    24 = @{ D = "0.001264";   E = "3.22%"  }
    25 = @{ D = "0.004392";   E = "1.68%"  }
    26 = @{ D = "0.0001141";  E = "-11.74%" }
    39 = @{ D = "0.02527";    E = "0.95%"  }
    40 = @{ D = "0.05259";    E = "1.32%"  }
    41 = @{ D = "0.007441";   E = "-3.34%" }
    42 = @{ D = "0.009613";   E = "4.77%"  }
    44 = @{ D = "0.002118";   E = "-2.14%" }
    45 = @{ D = "0.009975";   E = "10.95%" }
    46 = @{ D = "0.00006737"; E = "1.59%"  }
    48 = @{ D = "0.003002";   E = "-10.27%" }
}

# Rows whose Volume(1h) (E) changes but Price (D) stays the same
$volumeOnly = @{
    20 = "0.46%"
    21 = "3.49%"
    43 = "-1.49%"
    47 = "-0.12%"
    49 = "-7.72%"
    50 = "-0.12%"
    51 = "-0.12%"
}

foreach ($row in $priceAndVolume.Keys) {
    $vals = $priceAndVolume[$row]
    Set-TextValue $ws.Range("D$row") $vals.D
    Set-TextValue $ws.Range("E$row") $vals.E
}

foreach ($row in $volumeOnly.Keys) {
    Set-TextValue $ws.Range("E$row") $volumeOnly[$row]
}
